$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look like plain numbers need to be
# forced to Text format first, so Excel keeps the literal text (matching
# trailing/leading zeros) instead of coercing to a numeric value.
$textCells = @("D11", "D50", "D25", "D49", "D29", "D6", "D44", "D14", "D22", "D34", "D40", "D45", "D5", "D47", "D35", "D26", "D17")
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

# Column D (Price) updates
$ws.Range("D2").Value = "67.664.05"
$ws.Range("D3").Value = "3.786.38"
$ws.Range("D5").Value = "596.04"
$ws.Range("D6").Value = "166.75"
$ws.Range("D7").Value = "3.783.10"
$ws.Range("D11").Value = "6.33"
$ws.Range("D14").Value = "35.98"
$ws.Range("D15").Value = "4.421.45"
$ws.Range("D16").Value = "3.793.55"
$ws.Range("D17").Value = "18.50"
$ws.Range("D18").Value = "67.629.32"
$ws.Range("D22").Value = "458.92"
$ws.Range("D25").Value = "83.40"
$ws.Range("D26").Value = "12.07"
$ws.Range("D29").Value = "10.00"
$ws.Range("D30").Value = "3.936.63"
$ws.Range("D34").Value = "29.57"
$ws.Range("D35").Value = "0.999"
$ws.Range("D40").Value = "0.991"
$ws.Range("D44").Value = "48.00"
$ws.Range("D45").Value = "43.91"
$ws.Range("D47").Value = "150.56"
$ws.Range("D49").Value = "26.82"
$ws.Range("D50").Value = "388.64"

# Column E (Volume 1h) updates
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("E10").Value = "  -0.37%  "
$ws.Range("E11").Value = "  -1.12%  "
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("E13").Value = "  -2.31%  "
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("E17").Value = "  +3.53%  "
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("E21").Value = "  -7.63%  "
$ws.Range("E22").Value = "  -1.34%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").Value = "  +3.46%  "
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("E26").Value = "  +1.86%  "
$ws.Range("E27").Value = "  -2.87%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  -1.23%  "
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("E32").Value = "  +4.37%  "
$ws.Range("E33").Value = "  -1.08%  "
$ws.Range("E34").Value = "  -1.18%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  -0.72%  "
$ws.Range("E37").Value = "  -0.46%  "
$ws.Range("E38").Value = "  -2.47%  "
$ws.Range("E39").Value = "  -0.64%  "
$ws.Range("E40").Value = "  -1.04%  "
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("E45").Value = "  +0.81%  "
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("E47").Value = "  +3.18%  "
$ws.Range("E48").Value = "  -1.68%  "
$ws.Range("E49").Value = "  +6.04%  "
$ws.Range("E50").Value = "  -0.71%  "
$ws.Range("E51").Value = "  -4.42%  "
